# SI2019_Exposiciones.xlsx - "Add files via upload" commit
#
# The author filled in grading scores (columns G:M, plus the N "Total"
# SUM formula) for several rows, added review comments on a handful of
# cells, and repositioned/re-zoomed the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Grading scores for rows 15, 16, 17, 75, 76, 89 (cols G..M) plus the
#    N = SUM(G:M) total column.
# ---------------------------------------------------------------------
$rows = @(15, 16, 17, 75, 76, 89)
$scoreTable = @(
    @(6,  0, 8, 4, 4, 2, 3),
    @(6, 10, 9, 4, 4, 2, 3),
    @(6, 10, 9, 4, 4, 2, 3),
    @(7, 10, 9, 3, 4, 3, 3),
    @(7, 10, 9, 4, 4, 3, 3),
    @(8, 10, 10, 4, 4, 3, 3)
)

$cols = @("G", "H", "I", "J", "K", "L", "M")

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $vals = $scoreTable[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
    $ws.Range("N$row").Formula = "=SUM(G$($row):M$($row))"
}

# ---------------------------------------------------------------------
# 2. Review comments left by "Autor".
# ---------------------------------------------------------------------
$fullNote = "Autor:`n-No tiene estructura documentos TES`n-Falta Bibliografía`n-Está desarreglado"
$shortNote = "Autor:`n-No tiene estructura documentos TES`n-Falta Bibliografía"
$csharpNote = "Autor:`nEjercicio en C# en lugar de Visual Basic"

$ws.Range("G15").AddComment($fullNote)
$ws.Range("G16").AddComment($fullNote)
$ws.Range("G17").AddComment($fullNote)
$ws.Range("G75").AddComment($shortNote)
$ws.Range("J75").AddComment($csharpNote)
$ws.Range("G76").AddComment($shortNote)

# ---------------------------------------------------------------------
# 3. Sheet view: zoom to 90%, scroll the frozen pane down near the
#    bottom of the data and leave L89 selected.
# ---------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 90
$win = $excel.ActiveWindow
$win.ScrollRow = 77
$win.ScrollColumn = 1
$ws.Range("L89").Select()
